$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Proof-read / grammar fixes to the "EVENT SUMMARY" outcome texts (column B, rows 15-19)

# Row 15 (EVENT SUMMARY HEADER outcome): "Turns out" -> "It seems that"
$ws.Range("B15").Value = "It seems that the neighbouring kingdom was developing new technology to improve the quality of their equipment and was intending to share this information with other kingdoms once it was completed."

# Row 16 (ACTION 1 EVENT SUMMARY outcome) - text unchanged, rewritten for completeness
$ws.Range("B16").Value = "The neighbouring kingdom was outraged when they caught your spy, implying the amount of distrust you had. This had worsened the mutal ties between the kingdoms."

# Row 17 (ACTION 2 EVENT SUMMARY outcome): "thy" -> "they", added comma after "Although"
$ws.Range("B17").Value = "During the meeting, the military commander of the neighbouring kingdom assured you that they are not producing weapons with the aim of attacking other kingdoms. Although, you could sense a little displeasure in their tone due to your insistent probing."

# Row 18 (ACTION 3 EVENT SUMMARY outcome): "the ears of the king from the neighbouring kingdom" -> "the ears of their king"
$ws.Range("B18").Value = "The news of your soldiers raiding the neighbouring kingdom's weapons facility reached the ears of their king. He had seen this move as a call for war and decides to wage war against you."

# Row 19 (ACTION 4 EVENT SUMMARY outcome): "they have all decided" -> "they decided"
$ws.Range("B19").Value = "Word of your weapons facility had spread far and wide throughout the neighbouring kingdoms and they had all seen this as a threat. As such, they decided to make an alliance to wage war on you in fear of you getting too powerful."

# Update the active selection as recorded in the saved workbook
$ws.Range("O18").Select()
